# increase pipe finite elements to 4
$wb = $excel.ActiveWorkbook

$wsPipes = $wb.Worksheets.Item("Pipes")

# Nvol (column G) goes from 2 to 4 for every pipe (rows 2-40)
$wsPipes.Range("G2:G40").Value = 4

# The author ended their session with the Pipes sheet active and the
# newly-edited Nvol column selected (G2:G40, active cell G2), scrolled
# down so row 7 is at the top.
$wsPipes.Activate()
$wsPipes.Range("G2:G40").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
